# Update Handback report timestamps (simulating a regenerated report)
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for c73640ef...md row
$wsOverview.Range("G2").Value = "2016-08-21 11:07:49"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for c73640ef...md row
$wsZhCn.Range("H2").Value = "2016-08-21 11:07:45"
$wsZhCn.Range("K2").Value = "2016-08-21 11:08:07"

# de-de sheet: Correspond Handoff Datetime stays same value text (shared with Overview),
# Correspond Handback DateTime updated for c73640ef...md row
$wsDeDe.Range("H2").Value = "2016-08-21 11:07:49"
$wsDeDe.Range("K2").Value = "2016-08-21 11:08:13"
